$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the worksheet tab (from "SA-HW45.xpc" to "SA")
$ws.Name = "SA"

# Add a new row 16 mirroring the pattern of existing rows (same data as row 15
# except the index column increments and the label reuses the same category).
$ws.Cells.Item(15, 1).Copy()
$ws.Cells.Item(16, 1).PasteSpecial(-4122)
$ws.Cells.Item(16, 1).Value = 14

$ws.Cells.Item(16, 2).Value = $ws.Cells.Item(15, 2).Value2

for ($col = 3; $col -le 16; $col++) {
    $ws.Cells.Item(16, $col).Value = 1
}
